$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates: South Sudan gains the "*" annotation, Nigeria loses it ---
$ws.Range("B34").Value = "South Sudan*"
$ws.Range("B57").Value = "Nigeria"

# --- Swap the "resource-rich" shading between South Sudan (row 34) and Nigeria (row 57) ---
# South Sudan becomes shaded (copy the shaded format from a resource-rich row, e.g. row 17 - Congo)
$ws.Range("B17:J17").Copy()
$ws.Range("B34:J34").PasteSpecial(-4122)
# Nigeria becomes unshaded (copy the plain format from a non-resource-rich row, e.g. row 10 - South Africa)
$ws.Range("B10:J10").Copy()
$ws.Range("B57:J57").PasteSpecial(-4122)

# --- Recalculated aggregate/group statistics affected by the resource-rich reclassification ---
# Row 69
$ws.Range("C69").Value = 81.512351428571407
$ws.Range("D69").Value = 80.553692857142906
$ws.Range("E69").Value = 82.596598571428601
$ws.Range("F69").Value = 0.97614714285713999
$ws.Range("G69").Value = 69.872294285714304
$ws.Range("H69").Value = 64.014755714285698
$ws.Range("I69").Value = 76.2134342857143
$ws.Range("J69").Value = 0.83153285714285996
# Row 77
$ws.Range("C77").Value = 98.848484545454596
$ws.Range("D77").Value = 98.978205454545503
$ws.Range("E77").Value = 98.722103636363698
$ws.Range("F77").Value = 1.00259090909091
$ws.Range("G77").Value = 94.380635454545498
$ws.Range("H77").Value = 93.5228781818182
$ws.Range("I77").Value = 95.265992727272703
$ws.Range("J77").Value = 0.98175454545454999
# Row 80
$ws.Range("C80").Value = 73.430973750000007
$ws.Range("D80").Value = 70.409528750000007
$ws.Range("E80").Value = 76.697427500000003
$ws.Range("F80").Value = 0.89090000000000003
$ws.Range("G80").Value = 63.690757499999997
$ws.Range("H80").Value = 56.873626250000001
$ws.Range("I80").Value = 70.794944999999998
$ws.Range("J80").Value = 0.76314249999999995
# Row 82
$ws.Range("C82").Value = 79.576575952381006
$ws.Range("D82").Value = 77.129447142857202
$ws.Range("E82").Value = 82.301341666666701
$ws.Range("F82").Value = 0.92154190476191
$ws.Range("G82").Value = 68.062124999999995
$ws.Range("H82").Value = 61.997815952381004
$ws.Range("I82").Value = 74.578712619047593
$ws.Range("J82").Value = 0.80435547619047998
# Row 83
$ws.Range("E83").Value = 97.182109393939399
# Row 84
$ws.Range("C84").Value = 67.861473913043497
$ws.Range("D84").Value = 63.851693478260898
$ws.Range("E84").Value = 72.357340869565206
$ws.Range("F84").Value = 0.86088217391304001
$ws.Range("G84").Value = 54.949372173913098
$ws.Range("H84").Value = 46.670636521739098
$ws.Range("I84").Value = 64.018209999999996
$ws.Range("J84").Value = 0.70187304347825996
# Row 86
$ws.Range("C86").Value = 85.234460476190506
$ws.Range("D86").Value = 83.312081428571503
$ws.Range("E86").Value = 87.272696666666704
$ws.Range("F86").Value = 0.94830857142856995
$ws.Range("G86").Value = 74.241241428571399
$ws.Range("H86").Value = 68.710060952380999
$ws.Range("I86").Value = 79.9204890476191
$ws.Range("J86").Value = 0.84779809523810001
# Row 87
$ws.Range("C87").Value = 94.772266153846203
$ws.Range("D87").Value = 94.578426153846195
$ws.Range("E87").Value = 94.950995384615396
$ws.Range("F87").Value = 0.99491192307692
$ws.Range("G87").Value = 85.559614230769199
$ws.Range("H87").Value = 82.605646538461599
$ws.Range("I87").Value = 88.611429615384594
$ws.Range("J87").Value = 0.92505499999999996
# Row 89
$ws.Range("C89").Value = 98.818251666666697
$ws.Range("D89").Value = 98.929923333333406
$ws.Range("E89").Value = 98.715930555555602
$ws.Range("F89").Value = 1.00211583333333
$ws.Range("G89").Value = 95.951861944444403
$ws.Range("H89").Value = 95.305734722222198
$ws.Range("I89").Value = 96.611142222222199
$ws.Range("J89").Value = 0.98613166666666996
# Row 90
$ws.Range("C90").Value = 99.108101428571402
$ws.Range("D90").Value = 99.327102857142904
$ws.Range("E90").Value = 98.965710000000101
$ws.Range("F90").Value = 1.0036700000000001
$ws.Range("G90").Value = 97.481631818181896
$ws.Range("H90").Value = 97.137345909090897
$ws.Range("I90").Value = 97.744622727272699
$ws.Range("J90").Value = 0.99384454545455003
# Row 91
$ws.Range("H91").Value = 49.7276545161291
# Row 97
$ws.Range("C97").Value = 72.173169705882401
$ws.Range("D97").Value = 68.432221176470605
$ws.Range("E97").Value = 76.316578235294102
$ws.Range("F97").Value = 0.87812823529411999
$ws.Range("G97").Value = 60.930381764705899
$ws.Range("H97").Value = 53.7724282352941
$ws.Range("I97").Value = 68.646050000000002
$ws.Range("J97").Value = 0.75284735294117999
# Row 98
$ws.Range("C98").Value = 90.127426249999999
$ws.Range("D98").Value = 88.889005624999996
$ws.Range("E98").Value = 91.412044374999994
$ws.Range("F98").Value = 0.96518187499999997
$ws.Range("G98").Value = 79.61834125
$ws.Range("H98").Value = 75.608509999999995
$ws.Range("I98").Value = 83.711968749999997
$ws.Range("J98").Value = 0.88613687500000005
